$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 313, shifting existing rows 313:438 down to 314:439
$ws.Rows.Item(313).Insert()

# Populate the new row 313 with data
$ws.Cells.Item(313, 1).Value = 10
$ws.Cells.Item(313, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(313, 3).Value = "La Araucanía"
$ws.Cells.Item(313, 4).Value = 45229
$ws.Cells.Item(313, 5).Value = 9
$ws.Cells.Item(313, 6).Value = 100112039
$ws.Cells.Item(313, 7).Value = "Ciboulette"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 110
$ws.Cells.Item(313, 11).Value = 6000
$ws.Cells.Item(313, 12).Value = 6000
$ws.Cells.Item(313, 13).Value = 6000
$ws.Cells.Item(313, 14).Value = "$/docena de atados"
$ws.Cells.Item(313, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(313, 16).Value = 2000
$ws.Cells.Item(313, 17).Value = 3
$ws.Cells.Item(313, 18).Value = "Hortaliza"
